$d = $word.ActiveDocument

# The "hearing type" paragraph currently reads: "by <<hearingType>>"
# It should become just "<<hearingType>>" - i.e. the leading "by " text
# (and the run(s) that carry it) must be removed, leaving the
# "<<hearingType>>" merge field run(s) untouched.
$r = $d.Content
$found = $r.Find.Execute("by <<hearingType>>", $true, $false, $false, $false, $false,
                          $true, 1, $false, "", 0)

if ($found) {
    $sel = $word.Selection
    $sel.SetRange($r.Start, $r.Start + 3)
    $sel.Delete()
}
